$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
Write-Host $ws.Name
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
Write-Host $excel.ActiveWindow.ScrollRow
Write-Host $excel.ActiveWindow.ScrollColumn

